$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H37").Value = 864.1429000000001
$ws.Range("J37").Value = 864.1429000000001
$ws.Range("L37").Value = 2592.4287
$ws.Range("N37").Value = -2844.4287
$ws.Range("H86").Value = 42833.4
$ws.Range("I86").Value = 50096.094
$ws.Range("J86").Value = 4704.25
$ws.Range("K86").Value = 50096.094
$ws.Range("L86").Value = 4704.25
$ws.Range("M86").Value = -48973.094
$ws.Range("N86").Value = -6950.25
$ws.Range("H89").Value = 42833.4
$ws.Range("I89").Value = 50096.094
$ws.Range("J89").Value = 4704.25
$ws.Range("K89").Value = 250480.47
$ws.Range("L89").Value = 23521.25
$ws.Range("M89").Value = -244864.47
$ws.Range("N89").Value = -34753.25
$ws.Range("H106").Value = 2447.4
$ws.Range("I106").Value = 2447.4
$ws.Range("K106").Value = 2447.4
$ws.Range("M106").Value = -1816.4
$ws.Range("H138").Value = 3449.6365
$ws.Range("I138").Value = 3070.05
$ws.Range("J138").Value = 7245.5
$ws.Range("K138").Value = 9210.150000000001
$ws.Range("L138").Value = 21736.5
$ws.Range("M138").Value = -4070.150000000001
$ws.Range("N138").Value = -32016.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12080.424
$ws.Range("I32").Value = 9110.487999999999
$ws.Range("K32").Value = 9110.487999999999
$ws.Range("M32").Value = -8823.487999999999
$ws.Range("H61").Value = 3957.5488
$ws.Range("I61").Value = 3779.5571
$ws.Range("K61").Value = 3779.5571
$ws.Range("M61").Value = -3567.5571
$ws.Range("H136").Value = 3957.5488
$ws.Range("I136").Value = 3779.5571
$ws.Range("K136").Value = 11338.6713
$ws.Range("M136").Value = -8788.6713
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2166.7
$ws.Range("I20").Value = 2364.8572
$ws.Range("J20").Value = 1704.3334
$ws.Range("K20").Value = 2364.8572
$ws.Range("L20").Value = 1704.3334
$ws.Range("M20").Value = -2117.8572
$ws.Range("N20").Value = -2198.3334
$ws.Range("H134").Value = 2449.0227
$ws.Range("I134").Value = 2296.756
$ws.Range("K134").Value = 6890.268
$ws.Range("M134").Value = -4355.268
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6280.2104
$ws.Range("I31").Value = 4444
$ws.Range("K31").Value = 4444
$ws.Range("M31").Value = -4149
$ws.Range("H34").Value = 6280.2104
$ws.Range("I34").Value = 4444
$ws.Range("K34").Value = 4444
$ws.Range("M34").Value = -4242
$ws.Range("H69").Value = 29985.428
$ws.Range("I69").Value = 29985.428
$ws.Range("K69").Value = 29985.428
$ws.Range("M69").Value = -29236.428
$ws.Range("H72").Value = 29985.428
$ws.Range("I72").Value = 29985.428
$ws.Range("K72").Value = 89956.284
$ws.Range("M72").Value = -86212.284
$ws.Range("H105").Value = 4544
$ws.Range("I105").Value = 3649.3333
$ws.Range("K105").Value = 3649.3333
$ws.Range("M105").Value = -1902.3333
$ws.Range("H141").Value = 290552.8
$ws.Range("J141").Value = 290552.8
$ws.Range("L141").Value = 290552.8
$ws.Range("N141").Value = -300912.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2341.6
$ws.Range("I121").Value = 3174.75
$ws.Range("K121").Value = 9524.25
$ws.Range("M121").Value = -8214.25
$ws.Range("H132").Value = 1969.75
$ws.Range("J132").Value = 2093
$ws.Range("L132").Value = 18837
$ws.Range("N132").Value = -23897
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2745.8
$ws.Range("J80").Value = 3245.6667
$ws.Range("L80").Value = 3245.6667
$ws.Range("N80").Value = -5241.6667
$ws.Range("H83").Value = 2745.8
$ws.Range("J83").Value = 3245.6667
$ws.Range("L83").Value = 16228.3335
$ws.Range("N83").Value = -26212.3335
$ws.Range("H122").Value = 4099.077
$ws.Range("I122").Value = 3046.8333
$ws.Range("J122").Value = 4414.75
$ws.Range("K122").Value = 9140.499899999999
$ws.Range("L122").Value = 13244.25
$ws.Range("M122").Value = -6690.499899999999
$ws.Range("N122").Value = -18144.25
$ws.Range("H132").Value = 24754.34
$ws.Range("I132").Value = 26865.072
$ws.Range("J132").Value = 10331
$ws.Range("K132").Value = 80595.216
$ws.Range("L132").Value = 30993
$ws.Range("M132").Value = -78065.216
$ws.Range("N132").Value = -36053
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 486955.1
$ws.Range("I7").Value = 537518.8
$ws.Range("K7").Value = 537518.8
$ws.Range("M7").Value = -537406.8
$ws.Range("H34").Value = 335.5
$ws.Range("I34").Value = 335.5
$ws.Range("K34").Value = 335.5
$ws.Range("M34").Value = -163.5
$ws.Range("H40").Value = 5237.5557
$ws.Range("I40").Value = 5406.3335
$ws.Range("K40").Value = 5406.3335
$ws.Range("M40").Value = -5270.3335
$ws.Range("H122").Value = 4633.242
$ws.Range("I122").Value = 4271.357
$ws.Range("J122").Value = 4899.8945
$ws.Range("K122").Value = 12814.071
$ws.Range("L122").Value = 14699.6835
$ws.Range("M122").Value = -10364.071
$ws.Range("N122").Value = -19599.6835
$ws.Range("H126").Value = 486955.1
$ws.Range("I126").Value = 537518.8
$ws.Range("K126").Value = 1612556.4
$ws.Range("M126").Value = -1610086.4
$ws.Range("H132").Value = 78779.56
$ws.Range("I132").Value = 122447.3
$ws.Range("K132").Value = 367341.9
$ws.Range("M132").Value = -364811.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = ""
$ws.Range("H96").Value = 5920.857
$ws.Range("I96").Value = 4527.7144
$ws.Range("J96").Value = 7314
$ws.Range("K96").Value = 4527.7144
$ws.Range("L96").Value = 7314
$ws.Range("M96").Value = -3154.7144
$ws.Range("N96").Value = -10060
$ws.Range("H100").Value = 1196.091
$ws.Range("I100").Value = 1111.1428
$ws.Range("K100").Value = 2222.2856
$ws.Range("M100").Value = -1681.2856
$ws.Range("H122").Value = 4218.467
$ws.Range("I122").Value = 4234.75
$ws.Range("J122").Value = 4153.3335
$ws.Range("K122").Value = 12704.25
$ws.Range("L122").Value = 12460.0005
$ws.Range("M122").Value = -10254.25
$ws.Range("N122").Value = -17360.0005
$ws.Range("H126").Value = 52416.05
$ws.Range("I126").Value = 61033.59
$ws.Range("J126").Value = 3583.3333
$ws.Range("K126").Value = 183100.77
$ws.Range("L126").Value = 10749.9999
$ws.Range("M126").Value = -180630.77
$ws.Range("N126").Value = -15689.9999
$ws.Range("H132").Value = 32251.684
$ws.Range("I132").Value = 52407.09
$ws.Range("J132").Value = 8913.842000000001
$ws.Range("K132").Value = 157221.27
$ws.Range("L132").Value = 26741.526
$ws.Range("M132").Value = -154691.27
$ws.Range("N132").Value = -31801.526
$ws.Range("H136").Value = 4416.92
$ws.Range("I136").Value = 4046.7
$ws.Range("J136").Value = 5897.8
$ws.Range("K136").Value = 12140.1
$ws.Range("L136").Value = 17693.4
$ws.Range("M136").Value = -9590.099999999999
$ws.Range("N136").Value = -22793.4
